$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 20.23723175994358
$ws.Range("C2").Value = 9.452954855597929
$ws.Range("E2").Value = 14.37370296477407
$ws.Range("F2").Value = 40.91970542759623
$ws.Range("G2").Value = 42.46241287319923
$ws.Range("H2").Value = 17.44494420847082
$ws.Range("J2").Value = 8.357119370047588
$ws.Range("L2").Value = 12.40326378114587
$ws.Range("N2").Value = 18.66876304833223
$ws.Range("B3").Value = 19.80372883384151
$ws.Range("C3").Value = 9.156335475268996
$ws.Range("E3").Value = 14.39343620801631
$ws.Range("F3").Value = 40.88006059933464
$ws.Range("G3").Value = 42.30234282351739
$ws.Range("H3").Value = 17.47861192361294
$ws.Range("J3").Value = 8.37818912054585
$ws.Range("L3").Value = 12.38732978367635
$ws.Range("N3").Value = 18.74107769181861
$ws.Range("B4").Value = 19.53740505981618
$ws.Range("C4").Value = 8.967433267183049
$ws.Range("E4").Value = 14.40680180193981
$ws.Range("F4").Value = 40.86648542452524
$ws.Range("G4").Value = 42.21920803977655
$ws.Range("H4").Value = 17.50324745534015
$ws.Range("J4").Value = 8.391715303096193
$ws.Range("L4").Value = 12.37926835209663
$ws.Range("N4").Value = 18.78745510616503
$ws.Range("B5").Value = 19.42899897456628
$ws.Range("C5").Value = 8.88881474472864
$ws.Range("E5").Value = 14.41256314726464
$ws.Range("F5").Value = 40.86366257474394
$ws.Range("G5").Value = 42.18915738528109
$ws.Range("H5").Value = 17.51428024435513
$ws.Range("J5").Value = 8.397376067721451
$ws.Range("L5").Value = 12.37641888799811
$ws.Range("N5").Value = 18.80685291464309
$ws.Range("B6").Value = 19.41101043527
$ws.Range("C6").Value = 8.875663203447932
$ws.Range("E6").Value = 14.41353884260364
$ws.Range("F6").Value = 40.86335747507965
$ws.Range("G6").Value = 42.184399075298
$ws.Range("H6").Value = 17.51617214906753
$ws.Range("J6").Value = 8.398325033011906
$ws.Range("L6").Value = 12.37597211927367
$ws.Range("N6").Value = 18.81010407438154
$ws.Range("B7").Value = 19.53594234386314
$ws.Range("C7").Value = 8.966379537682661
$ws.Range("E7").Value = 14.40687822619976
$ws.Range("F7").Value = 40.86643638484026
$ws.Range("G7").Value = 42.21878724968741
$ws.Range("H7").Value = 17.50339222853808
$ws.Range("J7").Value = 8.391791043195843
$ws.Range("L7").Value = 12.37922815609661
$ws.Range("N7").Value = 18.78771469052056
$ws.Range("B8").Value = 20.08788669648107
$ws.Range("C8").Value = 9.3521238644769
$ws.Range("E8").Value = 14.38024806184002
$ws.Range("F8").Value = 40.90380289885913
$ws.Range("G8").Value = 42.40408930065492
$ws.Range("H8").Value = 17.45572847917068
$ws.Range("J8").Value = 8.364262254749036
$ws.Range("L8").Value = 12.39741393828194
$ws.Range("N8").Value = 18.6932881042384
$ws.Range("B9").Value = 21.16220242168693
$ws.Range("C9").Value = 10.05227375292161
$ws.Range("E9").Value = 14.33791226140816
$ws.Range("F9").Value = 41.06233973762988
$ws.Range("G9").Value = 42.88658160942133
$ws.Range("H9").Value = 17.39383419221116
$ws.Range("J9").Value = 8.314928154770731
$ws.Range("L9").Value = 12.44661804631618
$ws.Range("N9").Value = 18.52371388512607
$ws.Range("B10").Value = 21.93784272564972
$ws.Range("C10").Value = 10.5294002892041
$ws.Range("E10").Value = 14.3127988333959
$ws.Range("F10").Value = 41.23042971496498
$ws.Range("G10").Value = 43.31182876393977
$ws.Range("H10").Value = 17.3677679979987
$ws.Range("J10").Value = 8.281480373463394
$ws.Range("L10").Value = 12.49084789299408
$ws.Range("N10").Value = 18.40851875165318
$ws.Range("B11").Value = 22.28602367565734
$ws.Range("C11").Value = 10.73784990857374
$ws.Range("E11").Value = 14.30266715329063
$ws.Range("F11").Value = 41.31798034922466
$ws.Range("G11").Value = 43.52013440196104
$ws.Range("H11").Value = 17.36015249989701
$ws.Range("J11").Value = 8.266863801001168
$ws.Range("L11").Value = 12.51268147929115
$ws.Range("N11").Value = 18.35812769852649
$ws.Range("B12").Value = 22.41706867718748
$ws.Range("C12").Value = 10.81551160831713
$ws.Range("E12").Value = 14.29901576654846
$ws.Range("F12").Value = 41.3527132065157
$ws.Range("G12").Value = 43.60109946631201
$ws.Range("H12").Value = 17.35788053851347
$ws.Range("J12").Value = 8.261414424383883
$ws.Range("L12").Value = 12.52119167244826
$ws.Range("N12").Value = 18.33933334410698
$ws.Range("B13").Value = 22.38888355904685
$ws.Range("C13").Value = 10.79884296352434
$ws.Range("E13").Value = 14.29979392823647
$ws.Range("F13").Value = 41.34516288220861
$ws.Range("G13").Value = 43.5835703799973
$ws.Range("H13").Value = 17.35834260721489
$ws.Range("J13").Value = 8.262584246411276
$ws.Range("L13").Value = 12.51934814141096
$ws.Range("N13").Value = 18.34336827771152
$ws.Range("B14").Value = 22.2968216636902
$ws.Range("C14").Value = 10.74426490672641
$ws.Range("E14").Value = 14.30236304210475
$ws.Range("F14").Value = 41.32080628276915
$ws.Range("G14").Value = 43.5267539927545
$ws.Range("H14").Value = 17.3599533121697
$ws.Range("J14").Value = 8.266413764954313
$ws.Range("L14").Value = 12.51337678628718
$ws.Range("N14").Value = 18.35657572027874
$ws.Range("B15").Value = 22.2403226406204
$ws.Range("C15").Value = 10.71066737489201
$ws.Range("E15").Value = 14.30396080681818
$ws.Range("F15").Value = 41.30609236100235
$ws.Range("G15").Value = 43.49222210580486
$ws.Range("H15").Value = 17.36101964748112
$ws.Range("J15").Value = 8.26877058780237
$ws.Range("L15").Value = 12.50975058734159
$ws.Range("N15").Value = 18.36470307114647
$ws.Range("B16").Value = 21.91498289094695
$ws.Range("C16").Value = 10.51560131473691
$ws.Range("E16").Value = 14.31348692022974
$ws.Range("F16").Value = 41.22492993595412
$ws.Range("G16").Value = 43.29851012638816
$ws.Range("H16").Value = 17.36835121654013
$ws.Range("J16").Value = 8.282447608559204
$ws.Range("L16").Value = 12.48945514812634
$ws.Range("N16").Value = 18.41185226662352
$ws.Range("B17").Value = 21.71410599782217
$ws.Range("C17").Value = 10.39370451775441
$ws.Range("E17").Value = 14.31966150552764
$ws.Range("F17").Value = 41.17796893197306
$ws.Range("G17").Value = 43.18344403812178
$ws.Range("H17").Value = 17.37393682472317
$ws.Range("J17").Value = 8.290991054631561
$ws.Range("L17").Value = 12.47744052270674
$ws.Range("N17").Value = 18.44129086178354
$ws.Range("B18").Value = 21.59813556793215
$ws.Range("C18").Value = 10.32278537849022
$ws.Range("E18").Value = 14.32333466082555
$ws.Range("F18").Value = 41.15200266949365
$ws.Range("G18").Value = 43.11866309183436
$ws.Range("H18").Value = 17.37754871870473
$ws.Range("J18").Value = 8.295961430428154
$ws.Range("L18").Value = 12.47069154097843
$ws.Range("N18").Value = 18.45841262035307
$ws.Range("B19").Value = 21.5588000879981
$ws.Range("C19").Value = 10.29863593045738
$ws.Range("E19").Value = 14.32459924801068
$ws.Range("F19").Value = 41.14339073981794
$ws.Range("G19").Value = 43.09697169064398
$ws.Range("H19").Value = 17.37884014546229
$ws.Range("J19").Value = 8.297654019242399
$ws.Range("L19").Value = 12.46843431206979
$ws.Range("N19").Value = 18.46424234347674
$ws.Range("B20").Value = 21.73553532208403
$ws.Range("C20").Value = 10.40676447979216
$ws.Range("E20").Value = 14.31899161954681
$ws.Range("F20").Value = 41.18286001292136
$ws.Range("G20").Value = 43.19554827973219
$ws.Range("H20").Value = 17.37330089728983
$ws.Range("J20").Value = 8.290075755716463
$ws.Range("L20").Value = 12.47870281288861
$ws.Range("N20").Value = 18.43813747703274
$ws.Range("B21").Value = 22.32388530119599
$ws.Range("C21").Value = 10.76033064628123
$ws.Range("E21").Value = 14.3016034080346
$ws.Range("F21").Value = 41.3279176727276
$ws.Range("G21").Value = 43.5433862190551
$ws.Range("H21").Value = 17.35946359014034
$ws.Range("J21").Value = 8.26528662354249
$ws.Range("L21").Value = 12.51512417593746
$ws.Range("N21").Value = 18.35268858076354
$ws.Range("B22").Value = 22.70366658796093
$ws.Range("C22").Value = 10.98396860618866
$ws.Range("E22").Value = 14.29131876039607
$ws.Range("F22").Value = 41.43191836256759
$ws.Range("G22").Value = 43.78283930585853
$ws.Range("H22").Value = 17.35398714635214
$ws.Range("J22").Value = 8.249584201856127
$ws.Range("L22").Value = 12.54033783933028
$ws.Range("N22").Value = 18.29851863765863
$ws.Range("B23").Value = 22.50144558562508
$ws.Range("C23").Value = 10.86530056307268
$ws.Range("E23").Value = 14.29670929340089
$ws.Range("F23").Value = 41.37557516040273
$ws.Range("G23").Value = 43.65394817751353
$ws.Range("H23").Value = 17.35658311472075
$ws.Range("J23").Value = 8.257919423260107
$ws.Range("L23").Value = 12.52675320057948
$ws.Range("N23").Value = 18.3272773575902
$ws.Range("B24").Value = 21.72584862073896
$ws.Range("C24").Value = 10.40086268158011
$ws.Range("E24").Value = 14.31929409089233
$ws.Range("F24").Value = 41.1806455387978
$ws.Range("G24").Value = 43.19007167429593
$ws.Range("H24").Value = 17.3735871522101
$ws.Range("J24").Value = 8.290489379666575
$ws.Range("L24").Value = 12.47813163784892
$ws.Range("N24").Value = 18.43956250813239
$ws.Range("B25").Value = 20.87336312353847
$ws.Range("C25").Value = 9.869218407931221
$ws.Range("E25").Value = 14.34831055864329
$ws.Range("F25").Value = 40.91970542759623
$ws.Range("G25").Value = 42.74347105280243
$ws.Range("H25").Value = 17.40718049489537
$ws.Range("J25").Value = 8.327780340624894
$ws.Range("L25").Value = 12.43187450627878
$ws.Range("N25").Value = 18.56793064735884
